# Auto-generated update of cryptos list (price + volume %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.331.63"
$ws.Range("E2").Value = "  +1.50%  "
# Row 3
$ws.Range("D3").Value = "3.566.43"
$ws.Range("E3").Value = "  +2.16%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.71%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.44%  "
# Row 7
$ws.Range("D7").Value = "3.560.22"
$ws.Range("E7").Value = "  +1.99%  "
# Row 8
$ws.Range("E8").Value = "  +0.05%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "
# Row 10
$ws.Range("E10").Value = "  +5.64%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.40%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.97%  "
# Row 13
$ws.Range("E13").Value = "  +1.91%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.60%  "
# Row 15
$ws.Range("D15").Value = "4.165.25"
$ws.Range("E15").Value = "  +1.99%  "
# Row 16
$ws.Range("D16").Value = "3.558.26"
$ws.Range("E16").Value = "  +1.79%  "
# Row 17
$ws.Range("D17").Value = "68.193.27"
# Row 18
$ws.Range("E18").Value = "  +0.03%  "
# Row 19
$ws.Range("E19").Value = "  +6.35%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.64%  "
# Row 21
$ws.Range("E21").Value = "  +10.93%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.97%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.644"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.72%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.93%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000131"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.02%  "
# Row 26
$ws.Range("D26").Value = "3.703.37"
$ws.Range("E26").Value = "  +1.97%  "
# Row 28
$ws.Range("E28").Value = "  +4.45%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.01%  "
# Row 30
$ws.Range("E30").Value = "  +3.88%  "
# Row 31
$ws.Range("E31").Value = "  +9.36%  "
# Row 32
$ws.Range("E32").Value = "  +5.73%  "
# Row 33
$ws.Range("E33").Value = "  +0.15%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.86%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.05%  "
# Row 36
$ws.Range("E36").Value = "  +3.93%  "
# Row 37
$ws.Range("D37").Value = "3.556.90"
$ws.Range("E37").Value = "  +2.14%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.53%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.70%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.87%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.16%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.58%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.33%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.85%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.27%  "
# Row 48
$ws.Range("E48").Value = "  +6.58%  "
# Row 49
$ws.Range("E49").Value = "  +4.52%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.56%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.262"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.54%  "
